$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "END" was typed first (it is the lowest new shared-string index),
# then the rest of the column was filled in top-to-bottom order.
$ws.Cells.Item(14, 7).Value = "END"

# Values for the new column G (G1:G13)
$values = @{
    1  = "aaaa"
    2  = "bbb"
    3  = "ccc"
    4  = "ddd"
    5  = "fffff"
    6  = "gggg"
    7  = "h"
    8  = "ii"
    9  = "jjj"
    10 = "kk"
    11 = "LLL"
    12 = "mmm"
    13 = "nn"
}

foreach ($r in 1..13) {
    $ws.Cells.Item($r, 7).Value = $values[$r]
}

foreach ($r in 1..14) {
    if ($r -eq 2) { continue }
    $cell = $ws.Cells.Item($r, 7)
    $cell.Font.Name = "Tahoma"
    $cell.HorizontalAlignment = -4131
}

# G2 additionally gets a left border
$g2 = $ws.Range("G2")
$g2.Font.Name = "Tahoma"
$g2.Borders(7).LineStyle = 1
$g2.Borders(7).Weight = 2
$g2.HorizontalAlignment = -4131

$ws.Range("G1").Select()
